$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(63.63974129353218, 49.32, 76.34000000000003),
    @(6.634686567164172, 3.212000000000001, 10.93600000000001),
    @(12.52660696517412, 6.880000000000003, 19.12799999999998),
    @(9.057074626865669, 4.96, 14.52),
    @(33.24298507462672, 21.6, 46.78),
    @(58.27663681592036, 43.84399999999999, 72.63600000000005),
    @(78.49391044776112, 63.93600000000002, 89.10000000000002),
    @(29.85303482587054, 18.99199999999999, 41.29199999999998),
    @(79.08173134328348, 64.42000000000002, 89.35999999999996),
    @(72.18212935323373, 57.74800000000001, 82.69199999999995),
    @(54.82401990049746, 41.49199999999998, 68.46799999999999),
    @(79.53150248756212, 64.85599999999999, 89.28400000000005),
    @(76.92218905472618, 63.312, 87.03200000000001),
    @(11.93906467661692, 6.628000000000002, 19.22000000000001),
    @(33.38407960198992, 21.888, 46.676),
    @(54.83820895522376, 40.63599999999997, 68.54800000000002),
    @(12.52660696517412, 6.880000000000003, 19.12799999999998),
    @(73.90179104477606, 60.57999999999999, 84.16800000000001),
    @(38.013791044776, 25.636, 51.36000000000003),
    @(20.18354228855705, 11.464, 28.912),
    @(80.11315422885555, 64.66, 90.00400000000002),
    @(12.97596019900497, 7.352000000000004, 20.13199999999999),
    @(33.31737313432823, 21.684, 46.48800000000001),
    @(26.56366169154216, 15.73599999999999, 38.212),
    @(20.4151442786069, 12.06, 31.17999999999999),
    @(79.97701492537303, 64.37200000000003, 89.62),
    @(61.26429850746266, 46.82400000000004, 74.21199999999997),
    @(14.16234825870646, 8.016000000000004, 21.5),
    @(52.30334328358195, 39.65600000000003, 66.60000000000002),
    @(76.02696517412924, 62.85199999999996, 86.00800000000005),
    @(77.6799601990049, 63.028, 87.87200000000001),
    @(77.07012935323365, 63.30799999999996, 87.58),
    @(67.99810945273619, 54.70399999999999, 80.63599999999997),
    @(35.15795024875618, 23.22000000000001, 47.136),
    @(58.13703482587064, 42.6, 72.70000000000006),
    @(23.06523383084566, 13.96, 32.66399999999997),
    @(77.65357213930336, 63.91600000000003, 87.36),
    @(57.13140298507454, 43.26400000000002, 69.71199999999996),
    @(16.26089552238805, 9.443999999999996, 25.176),
    @(54.59395024875613, 40.544, 68.58799999999999),
    @(73.62441791044765, 60.364, 84.47200000000001),
    @(52.30334328358195, 39.65600000000003, 66.60000000000002),
    @(65.81327363184074, 50.74800000000004, 78.33200000000002),
    @(74.87755223880593, 62.46799999999999, 85.60399999999994),
    @(74.60742288557196, 60.33600000000005, 85.87600000000003),
    @(80.17679601990041, 64.08399999999999, 89.95200000000003),
    @(54.80927363184074, 41.42399999999998, 68.55199999999998),
    @(14.49066666666666, 8.640000000000006, 21.99599999999998),
    @(22.92268656716405, 14.112, 32.18399999999999),
    @(6.631741293532331, 3.204000000000001, 10.924),
    @(72.71584079601975, 58.26399999999997, 84.66800000000001),
    @(56.81512437810937, 42.11600000000002, 70.34800000000001),
    @(11.6070447761194, 6.456000000000002, 18.892),
    @(20.12527363184066, 11.94400000000001, 30.06399999999998),
    @(76.11864676616898, 62.79199999999998, 86.37599999999999),
    @(77.6482985074626, 63.11200000000001, 87.77200000000001),
    @(23.00921393034815, 13.408, 33.57199999999998),
    @(70.45834825870638, 56.37599999999996, 81.24799999999999),
    @(13.30889552238806, 7.448, 20.46799999999999),
    @(77.58013930348253, 63.93599999999999, 87.76399999999995),
    @(72.18212935323373, 57.74800000000001, 82.69199999999995),
    @(78.42792039800982, 64.748, 88.256),
    @(13.4537711442786, 7.688000000000001, 20.58799999999999),
    @(76.13265671641791, 60.42400000000002, 87.04799999999997),
    @(66.4406965174128, 52.172, 78.19200000000002),
    @(65.61351243781083, 52.04799999999999, 77.35599999999994),
    @(19.01870646766166, 10.536, 29.03999999999999),
    @(80.55542288557201, 64.60399999999996, 90.63600000000005),
    @(16.15834825870646, 10.184, 24.528),
    @(38.88161194029835, 25.748, 52.83200000000001),
    @(46.94463681592033, 32.25599999999997, 60.67199999999998)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $vals = $newValues[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
